# Refresh the cryptocurrency price list on Sheet1 (coinranking.com scrape)
# to the new snapshot.
#
# Every data row keeps column A (rank index) untouched; only Coin (B),
# Link (C), Price (D) and Volume/1h (E) text gets updated. Two coins
# (rows 19-20, Polkadot/BitcoinCash) also swapped rank order, so their
# whole B:E content is replaced rather than just the numbers.
#
# All of these cells hold plain TEXT in the workbook (prices such as
# "58.684.18" or "0.0₃0804" are not valid numbers, and the percent
# cells carry padding spaces), so every update below must land back as
# text too. Range.Value normally does this fine, but Excel's COM layer
# auto-detects "clean" decimal numbers (e.g. "6.26", "4.00", "0.0978")
# and silently coerces them to a Double, which mangles the text (drops
# trailing zeros, appends binary float noise, ...). For any replacement
# value that looks like a bare decimal number we prefix it with the
# standard Excel "force text" apostrophe before assigning so it stays
# text, matching the cell's original storage type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '58.803.37' },
    @{ Cell = 'E2'; Value = '  +0.23%  ' },
    @{ Cell = 'D3'; Value = '2.648.05' },
    @{ Cell = 'E3'; Value = '  +3.44%  ' },
    @{ Cell = 'E4'; Value = '  +0.31%  ' },
    @{ Cell = 'D5'; Value = '514.12' },
    @{ Cell = 'E5'; Value = '  +1.28%  ' },
    @{ Cell = 'D6'; Value = '143.66' },
    @{ Cell = 'E6'; Value = '  +0.10%  ' },
    @{ Cell = 'D7'; Value = '0.998' },
    @{ Cell = 'E7'; Value = '  -0.10%  ' },
    @{ Cell = 'D8'; Value = '0.566' },
    @{ Cell = 'E8'; Value = '  +1.81%  ' },
    @{ Cell = 'D9'; Value = '2.683.22' },
    @{ Cell = 'E9'; Value = '  +4.71%  ' },
    @{ Cell = 'D10'; Value = '6.26' },
    @{ Cell = 'E10'; Value = '  +0.84%  ' },
    @{ Cell = 'D11'; Value = '0.106' },
    @{ Cell = 'E11'; Value = '  +4.35%  ' },
    @{ Cell = 'E12'; Value = '  +1.16%  ' },
    @{ Cell = 'E13'; Value = '  -1.26%  ' },
    @{ Cell = 'D14'; Value = '3.115.95' },
    @{ Cell = 'E14'; Value = '  +3.71%  ' },
    @{ Cell = 'D15'; Value = '58.861.21' },
    @{ Cell = 'E15'; Value = '  +0.44%  ' },
    @{ Cell = 'D16'; Value = '20.93' },
    @{ Cell = 'E16'; Value = '  +1.67%  ' },
    @{ Cell = 'E17'; Value = '  +1.86%  ' },
    @{ Cell = 'D18'; Value = '2.679.93' },
    @{ Cell = 'E18'; Value = '  +4.83%  ' },
    @{ Cell = 'B19'; Value = 'BitcoinCash' },
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' },
    @{ Cell = 'D19'; Value = '346.39' },
    @{ Cell = 'E19'; Value = '  +4.03%  ' },
    @{ Cell = 'B20'; Value = 'Polkadot' },
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Cell = 'D20'; Value = '4.54' },
    @{ Cell = 'E20'; Value = '  +0.54%  ' },
    @{ Cell = 'D21'; Value = '10.41' },
    @{ Cell = 'E21'; Value = '  +3.50%  ' },
    @{ Cell = 'D22'; Value = '6.12' },
    @{ Cell = 'E22'; Value = '  +2.98%  ' },
    @{ Cell = 'D23'; Value = '0.999' },
    @{ Cell = 'E23'; Value = '  +0.25%  ' },
    @{ Cell = 'D24'; Value = '60.82' },
    @{ Cell = 'E24'; Value = '  +1.64%  ' },
    @{ Cell = 'D25'; Value = '0.419' },
    @{ Cell = 'E25'; Value = '  +2.93%  ' },
    @{ Cell = 'D26'; Value = '2.768.77' },
    @{ Cell = 'E26'; Value = '  +3.67%  ' },
    @{ Cell = 'D27'; Value = '0.994' },
    @{ Cell = 'E27'; Value = '  -0.45%  ' },
    @{ Cell = 'D28'; Value = '0.160' },
    @{ Cell = 'E28'; Value = '  +2.06%  ' },
    @{ Cell = 'D29'; Value = '0.0₃0808' },
    @{ Cell = 'E29'; Value = '  +3.62%  ' },
    @{ Cell = 'D30'; Value = '7.21' },
    @{ Cell = 'E30'; Value = '  +4.63%  ' },
    @{ Cell = 'D31'; Value = '0.999' },
    @{ Cell = 'E31'; Value = '  -0.05%  ' },
    @{ Cell = 'D32'; Value = '6.36' },
    @{ Cell = 'E32'; Value = '  +8.83%  ' },
    @{ Cell = 'D33'; Value = '18.88' },
    @{ Cell = 'E33'; Value = '  +1.91%  ' },
    @{ Cell = 'E34'; Value = '  +1.96%  ' },
    @{ Cell = 'D35'; Value = '149.55' },
    @{ Cell = 'E35'; Value = '  +0.14%  ' },
    @{ Cell = 'D36'; Value = '1.01' },
    @{ Cell = 'E36'; Value = '  +12.51%  ' },
    @{ Cell = 'D37'; Value = '4.00' },
    @{ Cell = 'E37'; Value = '  +1.44%  ' },
    @{ Cell = 'E38'; Value = '  +3.22%  ' },
    @{ Cell = 'D39'; Value = '36.64' },
    @{ Cell = 'E39'; Value = '  +1.93%  ' },
    @{ Cell = 'D40'; Value = '0.845' },
    @{ Cell = 'E40'; Value = '  +2.95%  ' },
    @{ Cell = 'D41'; Value = '3.67' },
    @{ Cell = 'E41'; Value = '  +4.33%  ' },
    @{ Cell = 'E42'; Value = '  +1.25%  ' },
    @{ Cell = 'D43'; Value = '0.619' },
    @{ Cell = 'E43'; Value = '  +2.00%  ' },
    @{ Cell = 'D44'; Value = '278.68' },
    @{ Cell = 'E44'; Value = '  -3.20%  ' },
    @{ Cell = 'D45'; Value = '0.994' },
    @{ Cell = 'E45'; Value = '  -0.37%  ' },
    @{ Cell = 'D46'; Value = '0.0978' },
    @{ Cell = 'E46'; Value = '  -0.28%  ' },
    @{ Cell = 'D47'; Value = '19.57' },
    @{ Cell = 'E47'; Value = '  +4.63%  ' },
    @{ Cell = 'D48'; Value = '0.0532' },
    @{ Cell = 'E48'; Value = '  -0.22%  ' },
    @{ Cell = 'E49'; Value = '  +1.19%  ' },
    @{ Cell = 'D50'; Value = '10.28' },
    @{ Cell = 'E50'; Value = '  -0.38%  ' },
    @{ Cell = 'D51'; Value = '1.996.44' },
    @{ Cell = 'E51'; Value = '  +4.44%  ' }
)

foreach ($u in $updates) {
    $value = $u.Value
    if ($value -match '^-?\d+(\.\d+)?$') {
        # Looks like a plain number to Excel's type-inference - keep it text.
        $value = "'" + $value
    }
    $ws.Range($u.Cell).Value = $value
}
